# Update "new Madigan bike hours" - Riders and Average columns
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Riders (column C) and Average (column D) new values for rows 2-8
$ws.Range("C2").Value = 306
$ws.Range("D2").Value = 277.33

$ws.Range("C3").Value = 212
$ws.Range("D3").Value = 227

$ws.Range("C4").Value = 264
$ws.Range("D4").Value = 255.33

$ws.Range("C5").Value = 261
$ws.Range("D5").Value = 249.33

$ws.Range("C6").Value = 246
$ws.Range("D6").Value = 253.67

$ws.Range("C7").Value = 127
$ws.Range("D7").Value = 116.5

$ws.Range("C8").Value = 86
$ws.Range("D8").Value = 77
